# Auto-generated edit script to update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to stay text so numeric-looking strings
# (e.g. "0.9997", "313.62") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.615.82"
$ws.Range("E2").Value = "  +2.11%  "

$ws.Range("D3").Value = "1.695.11"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.77%  "

$ws.Range("D5").Value = "313.62"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("D7").Value = "0.3952"
$ws.Range("E7").Value = "  +1.25%  "

$ws.Range("D8").Value = "0.4036"
$ws.Range("E8").Value = "  +2.59%  "

$ws.Range("D9").Value = "56.79"
$ws.Range("E9").Value = "  +16.70%  "

$ws.Range("D10").Value = "0.9992"
$ws.Range("E10").Value = "  -0.85%  "

$ws.Range("D11").Value = "1.519"
$ws.Range("E11").Value = "  +9.18%  "

$ws.Range("D12").Value = "0.08772"
$ws.Range("E12").Value = "  +1.80%  "

$ws.Range("D13").Value = "7.326"
$ws.Range("E13").Value = "  +13.74%  "

$ws.Range("D14").Value = "23.10"
$ws.Range("E14").Value = "  +2.14%  "

$ws.Range("D15").Value = "0.00001319"
$ws.Range("E15").Value = "  +2.07%  "

$ws.Range("D16").Value = "7.614"
$ws.Range("E16").Value = "  +7.26%  "

$ws.Range("D17").Value = "1.692.06"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").Value = "100.41"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").Value = "0.07060"
$ws.Range("E19").Value = "  +4.40%  "

$ws.Range("D20").Value = "19.48"
$ws.Range("E20").Value = "  +3.17%  "

$ws.Range("D21").Value = "6.717"
$ws.Range("E21").Value = "  +1.76%  "

$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").Value = "14.17"
$ws.Range("E23").Value = "  +4.14%  "

$ws.Range("D24").Value = "24.594.51"
$ws.Range("E24").Value = "  +2.10%  "

$ws.Range("D25").Value = "3.027"
$ws.Range("E25").Value = "  +12.27%  "

$ws.Range("D26").Value = "2.309"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").Value = "22.38"
$ws.Range("E27").Value = "  +3.22%  "

$ws.Range("E28").Value = "  +1.26%  "

$ws.Range("D29").Value = "5.183"
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").Value = "133.41"
$ws.Range("E30").Value = "  +3.00%  "

$ws.Range("D31").Value = "7.590"
$ws.Range("E31").Value = "  +34.72%  "

$ws.Range("D32").Value = "1.880.27"
$ws.Range("E32").Value = "  +1.67%  "

$ws.Range("D33").Value = "1.086"
$ws.Range("E33").Value = "  -3.35%  "

$ws.Range("D34").Value = "7.399"
$ws.Range("E34").Value = "  +21.20%  "

$ws.Range("D35").Value = "0.08578"
$ws.Range("E35").Value = "  +1.12%  "

$ws.Range("D36").Value = "1.959"
$ws.Range("E36").Value = "  +10.31%  "

$ws.Range("D37").Value = "11.02"
$ws.Range("E37").Value = "  +7.16%  "

$ws.Range("D38").Value = "0.2714"
$ws.Range("E38").Value = "  +4.35%  "

$ws.Range("D39").Value = "14.73"
$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("D40").Value = "0.02744"
$ws.Range("E40").Value = "  +10.66%  "

$ws.Range("D41").Value = "0.09000"
$ws.Range("E41").Value = "  +2.83%  "

$ws.Range("D42").Value = "1.472"
$ws.Range("E42").Value = "  +3.29%  "

$ws.Range("D43").Value = "0.7637"
$ws.Range("E43").Value = "  +4.53%  "

$ws.Range("D44").Value = "0.7166"
$ws.Range("E44").Value = "  +3.71%  "

$ws.Range("D45").Value = "15.31"
$ws.Range("E45").Value = "  +3.49%  "

$ws.Range("D46").Value = "2.453"
$ws.Range("E46").Value = "  +4.48%  "

$ws.Range("D47").Value = "4.169"
$ws.Range("E47").Value = "  +2.36%  "

$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("B49").Value = "Flow"
$ws.Range("C49").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D49").Value = "1.326"
$ws.Range("E49").Value = "  +18.32%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "140.49"
$ws.Range("E50").Value = "  +1.09%  "

$ws.Range("D51").Value = "0.00000000382"
$ws.Range("E51").Value = "  +2.76%  "
